$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated player roster: name, position, team for rows 2-19
$rows = @(
  @(2,  "Cade Cunningham",     "PG,SG",    "Detroit Pistons"),
  @(3,  "Derrick White",       "PG,SG",    "Boston Celtics"),
  @(4,  "Jaden McDaniels",     "SF,PF",    "Minnesota Timberwolves"),
  @(5,  "Ausar Thompson",      "SF,PF",    "Detroit Pistons"),
  @(6,  "Julius Randle",       "PF,C",     "Minnesota Timberwolves"),
  @(7,  "Naz Reid",            "PF,C",     "Minnesota Timberwolves"),
  @(8,  "Deandre Ayton",       "C",        "Portland Trail Blazers"),
  @(9,  "Cameron Johnson",     "SF,PF",    "Brooklyn Nets"),
  @(10, "LaMelo Ball",         "PG,SG",    "Charlotte Hornets"),
  @(11, "Anthony Davis",       "PF,C",     "Los Angeles Lakers"),
  @(12, "Damian Lillard",      "PG",       "Milwaukee Bucks"),
  @(13, "Collin Sexton",       "PG,SG",    "Utah Jazz"),
  @(14, "Cole Anthony",        "PG",       "Orlando Magic"),
  @(15, "Bam Adebayo",         "C",        "Miami Heat"),
  @(16, "Malik Monk",          "PG,SG,SF", "Sacramento Kings"),
  @(17, "Isaiah Hartenstein",  "C",        "Oklahoma City Thunder"),
  @(18, "Jerami Grant",        "SF,PF",    "Portland Trail Blazers"),
  @(19, "Brandon Miller",      "SG,SF,PF", "Charlotte Hornets")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
